$d = $word.ActiveDocument

# In the "Requisitos" bullet list, move the
# "LOQ4083 - Fenômenos de Transporte I (Requisito fraco)" line from the
# end of the list to the beginning (ahead of "LOB1006 - Cálculo IV").
#
# Each line in the list is its own run, followed by a manual line break
# (represented as ^l in Word's Find syntax). Locate the whole three-line
# block first.
$target = $d.Content
$found = $target.Find.Execute(
    "LOB1006 -  Cálculo IV  (Requisito fraco)^lLOB1019 -  Física II  (Requisito fraco)^lLOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)^l",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Remove the existing three runs, then rebuild them (one run per
    # line, each ending in a manual line break) in the new order, on the
    # now-collapsed (insertion-point) range so the original ListBullet
    # paragraph is preserved exactly with no leftover/duplicated runs.
    $target.Delete()

    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
              '<w:r><w:t>LOQ4083 -  Fen&#244;menos de Transporte I  (Requisito fraco)</w:t><w:br/></w:r>' +
              '<w:r><w:t>LOB1006 -  C&#225;lculo IV  (Requisito fraco)</w:t><w:br/></w:r>' +
              '<w:r><w:t>LOB1019 -  F&#237;sica II  (Requisito fraco)</w:t><w:br/></w:r>' +
              '</w:p>'

    $target.InsertXML($newXml)
}
